$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283, shifting all subsequent rows (old 283-397) down to 284-398
$ws.Rows("283:283").Insert()

# Populate the newly inserted row 283 with the new record's data
$ws.Range("A283").Value = 4
$ws.Range("B283").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C283").Value = "Los Lagos"
$ws.Range("D283").Value = 44755
$ws.Range("D283").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E283").Value = 10
$ws.Range("F283").Value = 100114013
$ws.Range("G283").Value = "Zanahoria"
$ws.Range("H283").Value = "Sin especificar"
$ws.Range("I283").Value = "Primera"
$ws.Range("J283").Value = 150
$ws.Range("K283").Value = 10000
$ws.Range("L283").Value = 10000
$ws.Range("M283").Value = 10000
$ws.Range("N283").Value = "$/saco 20 kilos"
$ws.Range("O283").Value = "Provincia de Llanquihue"
$ws.Range("P283").Value = 500
$ws.Range("Q283").Value = 20
$ws.Range("R283").Value = "Hortaliza"
